# Loan RBI, Variable Instalments
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# before column N, shifting the existing "Late" / "heading" / "Outstanding"
# columns (N,O,P) one place to the right (O,P,Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column in front of column N - this shifts the
# existing N/O/P columns (and all their data) one column to the right.
$ws.Columns("N").Insert()

# The freshly inserted column inherits the width of its left neighbour
# (column M) as an explicit custom width, matching what Excel does when a
# column is inserted next to a bestFit column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Restore the active selection to reflect where the user ended up working.
$ws.Range("J17").Select() | Out-Null
